# Fruta / hortaliza, semanal
# Apply the row-wise permutation of Fecha (D), Volumen (M) and the
# associated price columns (N, O, P, S) for rows 3-13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Row = 3;  D = 44476; M = 80  },
    @{ Row = 4;  D = 44417; M = 80  },
    @{ Row = 5;  D = 44405; M = 50;  N = 1200; O = 1200; P = 1200; S = 1200 },
    @{ Row = 6;  D = 44424; M = 50;  N = 1200; O = 1200; P = 1200; S = 1200 },
    @{ Row = 7;  D = 44435; M = 130; N = 1300; O = 1300; P = 1300; S = 1300 },
    @{ Row = 8;  D = 44438; M = 60;  N = 1200; O = 1200; P = 1200; S = 1200 },
    @{ Row = 9;  D = 44343; M = 60;  N = 1300; O = 1300; P = 1300; S = 1300 },
    @{ Row = 10; D = 44431; M = 100; N = 1300; O = 1300; P = 1300; S = 1300 },
    @{ Row = 11; D = 44473; M = 120 },
    @{ Row = 12; D = 44357; M = 35;  N = 1000; O = 1000; P = 1000; S = 1000 },
    @{ Row = 13; D = 44418; M = 40  }
)

foreach ($change in $changes) {
    $r = $change.Row
    $ws.Cells.Item($r, 4).Value = $change.D   # D = Fecha
    $ws.Cells.Item($r, 13).Value = $change.M  # M = Volumen

    if ($change.ContainsKey('N')) { $ws.Cells.Item($r, 14).Value = $change.N }  # N = Precio mínimo
    if ($change.ContainsKey('O')) { $ws.Cells.Item($r, 15).Value = $change.O }  # O = Precio máximo
    if ($change.ContainsKey('P')) { $ws.Cells.Item($r, 16).Value = $change.P }  # P = Precio promedio ponderado
    if ($change.ContainsKey('S')) { $ws.Cells.Item($r, 19).Value = $change.S }  # S = Precio $/Kg
}
